$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking price strings we write are stored as TEXT
# (matching the original inline-string / text cell type) rather than being
# auto-coerced into Excel numbers. NumberFormat must be set to "@" BEFORE
# the value is assigned for this to take effect.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin (already non-numeric text, no format change needed)
$ws.Range("D2").Value = "36.673.57"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.961.80"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.45%  "

# Row 5 - BNB
Set-TextValue "D5" "243.95"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6 - XRP
Set-TextValue "D6" "0.613"
$ws.Range("E6").Value = "  +0.77%  "

# Row 7 - Solana
Set-TextValue "D7" "57.93"
$ws.Range("E7").Value = "  +1.54%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.25%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.370"
$ws.Range("E9").Value = "  +1.11%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0808"
$ws.Range("E10").Value = "  -2.77%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.20%  "

# Row 12 - Avalanche
Set-TextValue "D12" "21.94"
$ws.Range("E12").Value = "  +2.46%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.252.89"
$ws.Range("E13").Value = "  +2.32%  "

# Row 14 - Polygon
Set-TextValue "D14" "0.816"
$ws.Range("E14").Value = "  -0.07%  "

# Row 15 - Chainlink
Set-TextValue "D15" "13.61"
$ws.Range("E15").Value = "  +1.80%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.23"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.959.73"
$ws.Range("E17").Value = "  +3.22%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "36.599.62"
$ws.Range("E18").Value = "  +0.94%  "

# Row 19 - Litecoin
Set-TextValue "D19" "69.55"
$ws.Range("E19").Value = "  +0.71%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("E20").Value = "  -0.26%  "

# Row 21 - Uniswap
Set-TextValue "D21" "5.06"
$ws.Range("E21").Value = "  +1.41%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "227.14"
$ws.Range("E22").Value = "  -0.56%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.02%  "

# Row 24 - was PancakeSwap, now Toncoin
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D24" "2.35"
$ws.Range("E24").Value = "  +3.00%  "

# Row 25 - was Toncoin, now PancakeSwap
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "2.38"
$ws.Range("E25").Value = "  -2.71%  "

# Row 26 - Cosmos
Set-TextValue "D26" "9.28"
$ws.Range("E26").Value = "  -0.71%  "

# Row 27 - was Kaspa, now Monero
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D27" "160.59"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28 - was Monero, now Kaspa
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D28" "0.136"
$ws.Range("E28").Value = "  +12.19%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "19.28"
$ws.Range("E29").Value = "  -0.34%  "

# Row 30 - Stellar
Set-TextValue "D30" "0.118"
$ws.Range("E30").Value = "  +0.68%  "

# Row 31 - ImmutableX
Set-TextValue "D31" "1.11"
$ws.Range("E31").Value = "  -2.37%  "

# Row 32 - Filecoin
Set-TextValue "D32" "4.65"
$ws.Range("E32").Value = "  -0.53%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0615"
$ws.Range("E33").Value = "  -2.11%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "4.22"
$ws.Range("E34").Value = "  -1.83%  "

# Row 35 - THORChain
Set-TextValue "D35" "6.24"
$ws.Range("E35").Value = "  +4.39%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  +0.36%  "

# Row 37 - RenderToken
Set-TextValue "D37" "3.38"
$ws.Range("E37").Value = "  +15.57%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  +2.43%  "

# Row 39 - WEMIXToken
Set-TextValue "D39" "1.77"
$ws.Range("E39").Value = "  -0.25%  "

# Row 40 - Cronos
Set-TextValue "D40" "0.0997"
$ws.Range("E40").Value = "  +3.12%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +1.69%  "

# Row 42 - VeChain
Set-TextValue "D42" "0.0211"
$ws.Range("E42").Value = "  +1.68%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -1.12%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "15.92"
$ws.Range("E44").Value = "  +1.19%  "

# Row 45 - ARBITRUM
$ws.Range("E45").Value = "  +0.44%  "

# Row 46 - Maker
$ws.Range("D46").Value = "1.346.97"
$ws.Range("E46").Value = "  +0.59%  "

# Row 47 - Aave
Set-TextValue "D47" "87.11"
$ws.Range("E47").Value = "  -0.62%  "

# Row 48 - FraxShare
Set-TextValue "D48" "7.09"
$ws.Range("E48").Value = "  -1.09%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  +0.93%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.144.13"
$ws.Range("E50").Value = "  +2.42%  "

# Row 51 - MultiversX
Set-TextValue "D51" "43.23"
$ws.Range("E51").Value = "  -5.08%  "
